# Apply the edit described by the diff:
#  - Rename column header B1 from "session_file_name" to "csv_file_name"
#  - Remove every data row whose session_file_date (column A) equals
#    "April 05, 2024" (these sessions fall outside the evaluation window),
#    shifting remaining rows up so the table stays contiguous
#    (dimension shrinks from A1:G52 to A1:G44)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the header text for column B
$ws.Range("B1").Value = "csv_file_name"

# Find the last used row in column A and walk upward, deleting any row
# whose date column matches the value being dropped. Going bottom-up
# means deleting a row never invalidates the row indices we still have
# left to check.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = $lastRow; $r -ge 2; $r--) {
    $dateValue = $ws.Cells.Item($r, 1).Value2
    if ($dateValue -eq "April 05, 2024") {
        $ws.Rows.Item($r).Delete()
    }
}
